# Add a new "16-jun" column (L) to the right of the existing "15-jun" column (K),
# mirroring the same header style and the per-row values, then update the
# active selection to reflect where the user ended up (N5).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell for the new day column
$ws.Range("L1").Value = "16-jun"

# New day's values, one per data row (rows 2-18)
$values = @{
    2  = 0
    3  = 14.945487314099447
    4  = 19.085022473555476
    5  = 17.572681940411279
    6  = 0
    7  = 5.3498893018430147
    8  = 2.8379637672919702
    9  = 10.753019815331831
    10 = 15.154081275730999
    11 = 13.284348698237924
    12 = 0
    13 = 15.116935131384595
    14 = 0
    15 = 0
    16 = 18.139681428802216
    17 = 0
    18 = 0
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 12).Value = $values[$row]
}

# Update the current selection like the author left it after editing
$ws.Range("N5").Select()
